$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated s-val data (filtering save games) for row 2
$ws.Range("B2").Value = 1.505614041169197
$ws.Range("C2").Value = 1.65323645889881
$ws.Range("D2").Value = 0.1529057820181812
$ws.Range("E2").Value = 0.4998867070740569
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 3.811642989160245
